# Revert "Powerpoint writer: consolidate text run nodes."
#
# Split each affected title/caption paragraph's runs back into one run
# per "word" and one run per run-of-whitespace, e.g.
#   "Slide " + "1 " + "(Content)"
# becomes
#   "Slide" + " " + "1" + " " + "(Content)"
# The run properties (<a:rPr/>) are empty throughout this deck, so a
# plain re-tokenization (without trying to carry over rich formatting)
# reproduces the target structure exactly.

$p = $ppt.ActivePresentation

# (SlideIndex, ShapeIndex, ParagraphIndex) for every paragraph whose runs
# contain an internal/trailing space and therefore get split below.
$targets = @(
    @{ Slide = 1;  Shape = 1; Para = 1 },
    @{ Slide = 2;  Shape = 1; Para = 1 },
    @{ Slide = 3;  Shape = 1; Para = 1 },
    @{ Slide = 4;  Shape = 1; Para = 1 },
    @{ Slide = 5;  Shape = 1; Para = 1 },
    @{ Slide = 6;  Shape = 1; Para = 1 },
    @{ Slide = 6;  Shape = 3; Para = 1 },  # "an image" textbox
    @{ Slide = 7;  Shape = 1; Para = 1 },
    @{ Slide = 7;  Shape = 4; Para = 1 },  # "An image" textbox
    @{ Slide = 8;  Shape = 1; Para = 1 },
    @{ Slide = 8;  Shape = 4; Para = 1 },  # "An image" textbox
    @{ Slide = 9;  Shape = 1; Para = 1 },
    @{ Slide = 10; Shape = 1; Para = 1 },
    @{ Slide = 11; Shape = 1; Para = 1 }
)

foreach ($t in $targets) {
    $tr = $p.Slides.Item($t.Slide).Shapes.Item($t.Shape).TextFrame.TextRange.Paragraphs($t.Para, 1)
    $text = $tr.Text

    # Tokenize into maximal runs of non-whitespace and maximal runs of
    # whitespace; re-assigning identical text to each sub-range forces
    # the writer to emit a dedicated <a:r> per token instead of leaving
    # the original (already-merged) run intact.
    $tokens = [regex]::Matches($text, '\S+|\s+')
    foreach ($tok in $tokens) {
        $start = $tok.Index + 1   # Characters() is 1-based
        $len = $tok.Length
        $tr.Characters($start, $len).Text = $tok.Value
    }
}
